$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'27.328.79"
$ws.Range("E2").Value = '  -4.27%  '

$ws.Range("D3").Formula = "'1.861.29"
$ws.Range("E3").Value = '  -5.16%  '

$ws.Range("E4").Value = '  -1.19%  '

$ws.Range("D5").Formula = "'323.58"
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("E6").Value = '  -0.98%  '

$ws.Range("D7").Formula = "'0.4518"
$ws.Range("E7").Value = '  -5.88%  '

$ws.Range("D8").Formula = "'0.3866"
$ws.Range("E8").Value = '  -5.08%  '

$ws.Range("D9").Formula = "'48.11"
$ws.Range("E9").Value = '  -11.33%  '

$ws.Range("D10").Formula = "'0.07904"
$ws.Range("E10").Value = '  -6.93%  '

$ws.Range("E11").Value = '  -3.90%  '

$ws.Range("D12").Formula = "'21.44"
$ws.Range("E12").Value = '  -4.44%  '

$ws.Range("D13").Formula = "'1.858.30"
$ws.Range("E13").Value = '  -2.45%  '

$ws.Range("D14").Formula = "'5.893"
$ws.Range("E14").Value = '  -4.79%  '

$ws.Range("E15").Value = '  -5.79%  '

$ws.Range("D16").Formula = "'1.0000"
$ws.Range("E16").Value = '  -1.24%  '

$ws.Range("D17").Formula = "'0.00001035"
$ws.Range("E17").Value = '  -3.55%  '

$ws.Range("D18").Formula = "'85.85"
$ws.Range("E18").Value = '  -5.67%  '

$ws.Range("D19").Formula = "'0.06526"
$ws.Range("E19").Value = '  -1.86%  '

$ws.Range("D20").Formula = "'17.20"
$ws.Range("E20").Value = '  -7.33%  '

$ws.Range("D21").Formula = "'0.9999"
$ws.Range("E21").Value = '  -0.95%  '

$ws.Range("D22").Formula = "'5.518"
$ws.Range("E22").Value = '  -6.00%  '

$ws.Range("D23").Formula = "'27.336.46"
$ws.Range("E23").Value = '  -4.19%  '

$ws.Range("D24").Formula = "'10.83"
$ws.Range("E24").Value = '  -5.85%  '

$ws.Range("D25").Formula = "'2.269"
$ws.Range("E25").Value = '  -1.42%  '

$ws.Range("D26").Formula = "'2.092.80"
$ws.Range("E26").Value = '  -2.44%  '

$ws.Range("D27").Formula = "'152.45"
$ws.Range("E27").Value = '  -2.37%  '

$ws.Range("D28").Formula = "'19.74"
$ws.Range("E28").Value = '  -2.97%  '

$ws.Range("E29").Value = '  -5.54%  '

$ws.Range("D30").Formula = "'5.504"
$ws.Range("E30").Value = '  -6.35%  '

$ws.Range("D31").Formula = "'120.54"
$ws.Range("E31").Value = '  -3.46%  '

$ws.Range("D32").Formula = "'1.497"
$ws.Range("E32").Value = '  +2.23%  '

$ws.Range("D33").Formula = "'0.09309"
$ws.Range("E33").Value = '  -3.88%  '

$ws.Range("D34").Formula = "'0.9373"
$ws.Range("E34").Value = '  -5.77%  '

$ws.Range("D35").Formula = "'3.603"
$ws.Range("E35").Value = '  -2.53%  '

$ws.Range("D36").Formula = "'5.283"
$ws.Range("E36").Value = '  -6.40%  '

$ws.Range("D37").Formula = "'0.02236"
$ws.Range("E37").Value = '  -4.50%  '

$ws.Range("D38").Formula = "'0.06005"
$ws.Range("E38").Value = '  -4.00%  '

$ws.Range("D39").Formula = "'1.223"
$ws.Range("E39").Value = '  -3.01%  '

$ws.Range("D40").Formula = "'8.254"
$ws.Range("E40").Value = '  -9.46%  '

$ws.Range("D41").Formula = "'0.9999"
$ws.Range("E41").Value = '  -0.93%  '

$ws.Range("D42").Formula = "'0.5916"
$ws.Range("E42").Value = '  -5.30%  '

$ws.Range("D43").Formula = "'0.1888"
$ws.Range("E43").Value = '  -1.61%  '

$ws.Range("D44").Formula = "'10.18"
$ws.Range("E44").Value = '  -9.51%  '

$ws.Range("D45").Formula = "'1.280"
$ws.Range("E45").Value = '  -5.65%  '

$ws.Range("D46").Formula = "'0.5641"
$ws.Range("E46").Value = '  -5.48%  '

$ws.Range("D47").Formula = "'12.06"
$ws.Range("E47").Value = '  -8.07%  '

$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Formula = "'3.369"
$ws.Range("E48").Value = '  -1.31%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Formula = "'1.927"
$ws.Range("E49").Value = '  -6.94%  '

$ws.Range("D50").Formula = "'0.06795"
$ws.Range("E50").Value = '  -0.72%  '

$ws.Range("D51").Formula = "'107.90"
$ws.Range("E51").Value = '  -3.10%  '
